# Generate Report for Handoff
# Renames the two tracked files, flips their status from "handed back" to
# "ready for handoff" (i.e. a fresh handoff just happened and no handback
# has come in yet), and refreshes the handoff timestamps/hashes.

$wb = $excel.ActiveWorkbook

$oldFile1 = "c43b3f1f-b159-4eef-8940-5e41d1f38fc6.md"
$newFile1 = "cb8b6387-ef20-4e9d-963c-e493c86ea27a.md"
$oldFile2 = "e2944927-23b0-4647-9543-ab3703d28b1e.md"
$newFile2 = "ffff9ba7ca40-387a-4f0e-8e98-85c9db8fb613.md"

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

$oldHoDate = "2016-08-26 09:07:11"
$newHoDate = "2016-08-26 09:08:30"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $newFile1
$ov.Range("B2").Value = "e2e\" + $newFile1
$ov.Range("A3").Value = $newFile2
$ov.Range("B3").Value = "e2e\" + $newFile2

$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus

$ov.Range("G2").Value = $newHoDate
$ov.Range("G3").Value = $newHoDate

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/055865d5c14721d68676b5a42ed469dbc1bd4b63/e2e/c43b3f1f-b159-4eef-8940-5e41d1f38fc6.md", "", "", "e2e\" + $newFile1)
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/055865d5c14721d68676b5a42ed469dbc1bd4b63/e2e/e2944927-23b0-4647-9543-ab3703d28b1e.md", "", "", "e2e\" + $newFile2)

$ov.Columns.Item(5).ColumnWidth = 16.38265482584637
$ov.Columns.Item(6).ColumnWidth = 16.38265482584637

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $newFile1
$zh.Range("C2").Value = $newStatus
$zh.Range("G2").Value = "cb8b6387-ef20-4e9d-963c-e493c86ea27a.379c70dac170965e32cac0c420ccc280e6549c9b.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-26 09:08:25"
$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

$zh.Range("A3").Value = $newFile2
$zh.Range("C3").Value = $newStatus
$zh.Range("F3").Value = "True"
$zh.Range("G3").Value = $zh.Range("G2").Value()
$zh.Range("H3").Value = $zh.Range("H2").Value()
$zh.Range("I3").Value = ""
$zh.Range("I3").Style = "Normal"
$zh.Range("J3").Value = ""
$zh.Range("K3").Value = "0001-01-01 00:00:00"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/055865d5c14721d68676b5a42ed469dbc1bd4b63/e2e/c43b3f1f-b159-4eef-8940-5e41d1f38fc6.md", "", "", $newFile1)
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/055865d5c14721d68676b5a42ed469dbc1bd4b63/e2e/e2944927-23b0-4647-9543-ab3703d28b1e.md", "", "", $newFile2)

$zh.Columns.Item(3).ColumnWidth = 16.38265482584637
$zh.Columns.Item(9).ColumnWidth = 17.817272004627068
$zh.Columns.Item(10).ColumnWidth = 20.872143700009268

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $newFile1
$de.Range("C2").Value = $newStatus
$de.Range("G2").Value = "cb8b6387-ef20-4e9d-963c-e493c86ea27a.379c70dac170965e32cac0c420ccc280e6549c9b.de-de.xlf"
$de.Range("H2").Value = $newHoDate
$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"

$de.Range("A3").Value = $newFile2
$de.Range("C3").Value = $newStatus
$de.Range("F3").Value = "True"
$de.Range("G3").Value = $de.Range("G2").Value()
$de.Range("H3").Value = $de.Range("H2").Value()
$de.Range("I3").Value = ""
$de.Range("I3").Style = "Normal"
$de.Range("J3").Value = ""
$de.Range("K3").Value = "0001-01-01 00:00:00"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/055865d5c14721d68676b5a42ed469dbc1bd4b63/e2e/c43b3f1f-b159-4eef-8940-5e41d1f38fc6.md", "", "", $newFile1)
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/055865d5c14721d68676b5a42ed469dbc1bd4b63/e2e/e2944927-23b0-4647-9543-ab3703d28b1e.md", "", "", $newFile2)

$de.Columns.Item(3).ColumnWidth = 16.38265482584637
$de.Columns.Item(9).ColumnWidth = 17.817272004627068
$de.Columns.Item(10).ColumnWidth = 20.872143700009268
